$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "before E2:" $ws.Range("E2").Value()
$ws.Range("E2").Value = "H"
Write-Host "after E2:" $ws.Range("E2").Value()

Write-Host "before G206:" $ws.Range("G206").Value()
$ws.Range("G206").Value = "Fernández"
Write-Host "after G206:" $ws.Range("G206").Value()
